$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Dose multiplier interval" column is being inserted at G, pushing
# the existing "Dose vol." (G) and "Force delay" (H) columns one place to
# the right, into H and I. Only the header row and the two data rows are
# populated, so shift those cells explicitly (right to left) instead of a
# full-sheet column insert.
$newHeader = "Dose multiplier interval"
$newValues = @{ 2 = 0.1; 3 = 0.05 }

for ($r = 1; $r -le 3; $r++) {
    $gCell = $ws.Cells.Item($r, 7)
    $hCell = $ws.Cells.Item($r, 8)
    $iCell = $ws.Cells.Item($r, 9)

    # Give the newly-used I column the same formatting G/H already use.
    $gCell.Copy()
    $iCell.PasteSpecial(-4122)  # xlPasteFormats

    $hVal = $hCell.Value2
    $gVal = $gCell.Value2

    $iCell.Value2 = $hVal
    $hCell.Value2 = $gVal
}

# Header + values for the new column G
$ws.Cells.Item(1, 7).Value2 = $newHeader
$ws.Cells.Item(2, 7).Value2 = $newValues[2]
$ws.Cells.Item(3, 7).Value2 = $newValues[3]

$excel.CutCopyMode = $false

# Update selection to reflect new focus on column G
$ws.Range("G1:G3").Select()
